# Update Daily Report: 2026-02-27
# Appends the new reporting day (date serial 46079) to Daily_Data,
# and refreshes the roll-up figures on Today_Summary and Monthly_Stats.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Daily_Data: append rows 314-337 for the new day.
#    Columns: A Date, B Region_Type, C PREV_TOTAL, D RECEIVED,
#             E WITHDRAWN, F NET_CHANGE, G ADJUSTMENT, H TOTAL_TODAY
# ---------------------------------------------------------------------
$wsDaily = $wb.Worksheets.Item("Daily_Data")

$dailyRows = @(
    @(314, 46079, "ASAHI DEPOSITORY LLC Registered", 23291615.992, 824168.84, 0, 824168.84, 0, 24115784.832),
    @(315, 46079, "ASAHI DEPOSITORY LLC Eligible", 1703429.248, 0, 0, 0, 0, 1703429.248),
    @(316, 46079, "BRINK'S, INC. Registered", 14477133.393, 0, 0, 0, 111428.054, 14588561.447),
    @(317, 46079, "BRINK'S, INC. Eligible", 40642521.76, 0, 0, 0, -111428.054, 40531093.706),
    @(318, 46079, "CNT DEPOSITORY, INC. Registered", 12170205.469, 0, 0, 0, 0, 12170205.469),
    @(319, 46079, "CNT DEPOSITORY, INC. Eligible", 13861333.923, 0, 485175.71, -485175.71, 0, 13376158.213),
    @(320, 46079, "DELAWARE DEPOSITORY Registered", 1532776.423, 0, 0, 0, 85768.116, 1618544.539),
    @(321, 46079, "DELAWARE DEPOSITORY Eligible", 15758772.402, 26231.01, 9999.81, 16231.2, -85768.116, 15689235.486),
    @(322, 46079, "HSBC BANK, USA Registered", 3387219.03, 0, 0, 0, 0, 3387219.03),
    @(323, 46079, "HSBC BANK, USA Eligible", 17764657.013, 0, 0, 0, 0, 17764657.013),
    @(324, 46079, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 273789.87, 0, 0, 0, 0, 273789.87),
    @(325, 46079, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 3295246.644, 0, 0, 0, 0, 3295246.644),
    @(326, 46079, "JP MORGAN CHASE BANK NA Registered", 12000343.77, 0, 0, 0, 0, 12000343.77),
    @(327, 46079, "JP MORGAN CHASE BANK NA Eligible", 142491363.683, 0, 647661.5, -647661.5, 0, 141843702.183),
    @(328, 46079, "LOOMIS INTERNATIONAL (US) LLC Registered", 6842629.447, 0, 0, 0, 964465.51, 7807094.957),
    @(329, 46079, "LOOMIS INTERNATIONAL (US) LLC Eligible", 23512931.636, 0, 0, 0, -964465.51, 22548466.126),
    @(330, 46079, "MALCA-AMIT ARMORED, INC. Registered", 0, 0, 0, 0, 0, 0),
    @(331, 46079, "MALCA-AMIT ARMORED, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @(332, 46079, "MALCA-AMIT USA, LLC Registered", 949634.064, 0, 0, 0, 0, 949634.064),
    @(333, 46079, "MALCA-AMIT USA, LLC Eligible", 1073898.377, 0, 128731, -128731, 0, 945167.377),
    @(334, 46079, "MANFRA, TORDELLA & BROOKES, LLC Registered", 5871594.333, 0, 0, 0, 0, 5871594.333),
    @(335, 46079, "MANFRA, TORDELLA & BROOKES, LLC Eligible", 11984893.898, 0, 192764.302, -192764.302, 0, 11792129.596),
    @(336, 46079, "STONEX PRECIOUS METALS LLC Registered", 5333457.37, 307500.19, 0, 307500.19, 0, 5640957.56),
    @(337, 46079, "STONEX PRECIOUS METALS LLC Eligible", 2419487.69, 0, 0, 0, 0, 2419487.69)
)

# Match the existing date-stamp formatting used by column A (style carries
# numFmtId 165 "YYYY-MM-DD HH:MM:SS") before writing the serial values.
$wsDaily.Range("A314:A337").NumberFormat = "YYYY-MM-DD HH:MM:SS"

foreach ($row in $dailyRows) {
    $r = $row[0]
    $wsDaily.Cells.Item($r, 1).Value = $row[1]
    $wsDaily.Cells.Item($r, 2).Value = $row[2]
    $wsDaily.Cells.Item($r, 3).Value = $row[3]
    $wsDaily.Cells.Item($r, 4).Value = $row[4]
    $wsDaily.Cells.Item($r, 5).Value = $row[5]
    $wsDaily.Cells.Item($r, 6).Value = $row[6]
    $wsDaily.Cells.Item($r, 7).Value = $row[7]
    $wsDaily.Cells.Item($r, 8).Value = $row[8]
}

# ---------------------------------------------------------------------
# 2) Today_Summary: refresh Eligible / Registered / Total_Stock per
#    depository using the new day's closing TOTAL_TODAY figures.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Today_Summary")

$summaryRows = @(
    @(2, 1703429.248, 24115784.832),
    @(3, 40531093.706, 14588561.447),
    @(4, 13376158.213, 12170205.469),
    @(5, 15689235.486, 1618544.539),
    @(6, 17764657.013, 3387219.03),
    @(7, 3295246.644, 273789.87),
    @(8, 141843702.183, 12000343.77),
    @(9, 22548466.126, 7807094.957),
    @(10, 0, 0),
    @(11, 945167.377, 949634.064),
    @(12, 11792129.596, 5871594.333),
    @(13, 2419487.69, 5640957.56)
)

foreach ($row in $summaryRows) {
    $r = $row[0]
    $eligible = $row[1]
    $registered = $row[2]
    $wsSummary.Cells.Item($r, 2).Value = $eligible
    $wsSummary.Cells.Item($r, 3).Value = $registered
    $wsSummary.Cells.Item($r, 4).Value = $eligible + $registered
}

# ---------------------------------------------------------------------
# 3) Monthly_Stats: refresh the month-to-date roll-up (row 2) and the
#    per-depository/type detail rows affected by the new day's activity.
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")

$wsMonthly.Cells.Item(2, 2).Value = 271908773.282
$wsMonthly.Cells.Item(2, 3).Value = 88423729.871
$wsMonthly.Cells.Item(2, 4).Value = 360332503.153

$monthlyRows = @(
    @(8, 824168.84, 0, 24115784.832),
    @(9, 0, 2457230.712, 40531093.706),
    @(10, 0, 0, 14588561.447),
    @(11, 7944.85, 5631138.443, 13376158.213),
    @(13, 275185.199, 714523.477, 15689235.486),
    @(14, 0, 0, 1618544.539),
    @(19, 0, 15020099.25, 141843702.183),
    @(21, 10089.96, 2964461.58, 22548466.126),
    @(22, 0, 0, 7807094.957),
    @(25, 0, 128731, 945167.377),
    @(27, 0, 1854483.234, 11792129.596),
    @(30, 392565.06, 0, 5640957.56)
)

foreach ($row in $monthlyRows) {
    $r = $row[0]
    $wsMonthly.Cells.Item($r, 3).Value = $row[1]
    $wsMonthly.Cells.Item($r, 4).Value = $row[2]
    $wsMonthly.Cells.Item($r, 5).Value = $row[3]
}
